# 24/04/2025 - Dang code them chuc nang comment
# Add a bold, red "reviewer comment" note in cell C12 of the "Special" sheet
# and refresh the sheet view scroll/selection position.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Special")

# New comment cell next to the "2. Update Task" description
$c = $ws.Range("C12")
$c.Value = "GỬI THÔNG BÁO CHO AI?"
$c.Font.Bold = $true
$c.Font.Color = 255
$c.HorizontalAlignment = -4131
$c.WrapText = $false

# Update the view's scroll position / active selection as seen when the
# edit was made
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 8
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("D14").Select()
